$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D and E are treated as text so values like "1.00" or "87.543.16" are not
# reinterpreted as numbers/dates by Excel when assigned via .Value
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '87.543.16'
$ws.Range('E2').Value = '  +1.28%  '
$ws.Range('D3').Value = '3.165.86'
$ws.Range('E3').Value = '  -3.33%  '
$ws.Range('E4').Value = '  +0.38%  '
$ws.Range('D5').Value = '207.89'
$ws.Range('E5').Value = '  -2.36%  '
$ws.Range('D6').Value = '607.10'
$ws.Range('E6').Value = '  -2.90%  '
$ws.Range('D7').Value = '0.390'
$ws.Range('E7').Value = '  +3.41%  '
$ws.Range('D8').Value = '0.675'
$ws.Range('E8').Value = '  +7.46%  '
$ws.Range('E9').Value = '  +0.09%  '
$ws.Range('D10').Value = '3.163.08'
$ws.Range('E10').Value = '  -3.25%  '
$ws.Range('D11').Value = '0.538'
$ws.Range('E11').Value = '  -7.05%  '
$ws.Range('E12').Value = '  +1.07%  '
$ws.Range('D13').Value = '0.0000244'
$ws.Range('E13').Value = '  -6.79%  '
$ws.Range('B14').Value = 'Toncoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D14').Value = '5.27'
$ws.Range('E14').Value = '  +1.07%  '
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '3.737.99'
$ws.Range('E15').Value = '  -3.32%  '
$ws.Range('D16').Value = '87.242.90'
$ws.Range('E16').Value = '  +1.19%  '
$ws.Range('D17').Value = '32.13'
$ws.Range('E17').Value = '  -6.48%  '
$ws.Range('D18').Value = '3.151.93'
$ws.Range('E18').Value = '  -4.17%  '
$ws.Range('D19').Value = '3.18'
$ws.Range('E19').Value = '  +9.02%  '
$ws.Range('D20').Value = '13.42'
$ws.Range('E20').Value = '  -5.03%  '
$ws.Range('D21').Value = '412.69'
$ws.Range('E21').Value = '  -4.83%  '
$ws.Range('D22').Value = '8.48'
$ws.Range('E22').Value = '  -7.23%  '
$ws.Range('D23').Value = '5.06'
$ws.Range('E23').Value = '  -5.19%  '
$ws.Range('D24').Value = '5.17'
$ws.Range('E24').Value = '  +0.04%  '
$ws.Range('D25').Value = '11.99'
$ws.Range('E25').Value = '  -0.29%  '
$ws.Range('D26').Value = '3.324.05'
$ws.Range('E26').Value = '  -4.30%  '
$ws.Range('B27').Value = 'PEPE'
$ws.Range('C27').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D27').Value = '0.0000132'
$ws.Range('E27').Value = '  +1.38%  '
$ws.Range('B28').Value = 'Litecoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D28').Value = '73.26'
$ws.Range('E28').Value = '  -4.38%  '
$ws.Range('E29').Value = '  +0.20%  '
$ws.Range('D30').Value = '0.161'
$ws.Range('E30').Value = '  -2.96%  '
$ws.Range('D31').Value = '0.996'
$ws.Range('E31').Value = '  -0.61%  '
$ws.Range('D32').Value = '546.66'
$ws.Range('E32').Value = '  +0.25%  '
$ws.Range('D33').Value = '8.22'
$ws.Range('E33').Value = '  -7.31%  '
$ws.Range('D34').Value = '1.31'
$ws.Range('E34').Value = '  -8.88%  '
$ws.Range('B35').Value = 'RenderToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D35').Value = '6.82'
$ws.Range('E35').Value = '  +0.38%  '
$ws.Range('B36').Value = 'PancakeSwap'
$ws.Range('C36').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D36').Value = '1.85'
$ws.Range('E36').Value = '  -6.73%  '
$ws.Range('D37').Value = '0.131'
$ws.Range('E37').Value = '  -3.65%  '
$ws.Range('D38').Value = '21.78'
$ws.Range('E38').Value = '  -3.39%  '
$ws.Range('D39').Value = '21.83'
$ws.Range('E39').Value = '  +0.41%  '
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  +0.83%  '
$ws.Range('D41').Value = '3.03'
$ws.Range('E41').Value = '  +2.13%  '
$ws.Range('E42').Value = '  +0.25%  '
$ws.Range('D43').Value = '1.92'
$ws.Range('E43').Value = '  -4.44%  '
$ws.Range('D44').Value = '0.369'
$ws.Range('E44').Value = '  -8.67%  '
$ws.Range('D45').Value = '149.72'
$ws.Range('E45').Value = '  -1.85%  '
$ws.Range('D46').Value = '173.12'
$ws.Range('E46').Value = '  -3.64%  '
$ws.Range('D47').Value = '43.09'
$ws.Range('E47').Value = '  -2.95%  '
$ws.Range('D48').Value = '0.125'
$ws.Range('E48').Value = '  +5.00%  '
$ws.Range('D49').Value = '1.23'
$ws.Range('E49').Value = '  -8.17%  '
$ws.Range('D50').Value = '3.97'
$ws.Range('E50').Value = '  -6.37%  '
$ws.Range('D51').Value = '0.694'
$ws.Range('E51').Value = '  -5.59%  '
